$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(180, 1).Value = 177
$ws.Cells.Item(180, 2).Value = 626402
$ws.Cells.Item(180, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(180, 4).Value = "2024-02-06"
$ws.Cells.Item(180, 5).Value = "''00063T"
$ws.Cells.Item(180, 6).Value = "''241751303001087"
$ws.Cells.Item(180, 7).Value = "Pembayaran belanja barangberupa biaya perjadin sesuai STNo.B.309,B.356,B.729/BPPSDM.1/KP.440/I/2024 Tanggal 15,16 dan 19 Januari 2024 a.n Achmad Irfansyah, dkk"
$ws.Cells.Item(180, 8).Value = "''626402.175.524111.03212WA.2378EBA.A000000001.00000.1.0151.2.000000.000000.957.101.0A.000311"
$ws.Cells.Item(180, 9).Value = "IDR"
$ws.Cells.Item(180, 10).Value = 1
$ws.Cells.Item(180, 11).Value = 2316500
$ws.Cells.Item(180, 12).Value = 2316500

$ws.Cells.Item(181, 1).Value = 178
$ws.Cells.Item(181, 2).Value = 626402
$ws.Cells.Item(181, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(181, 4).Value = "2024-02-06"
$ws.Cells.Item(181, 5).Value = "''00063T"
$ws.Cells.Item(181, 6).Value = "''241751303001087"
$ws.Cells.Item(181, 7).Value = "Pembayaran belanja barangberupa biaya perjadin sesuai STNo.B.309,B.356,B.729/BPPSDM.1/KP.440/I/2024 Tanggal 15,16 dan 19 Januari 2024 a.n Achmad Irfansyah, dkk"
$ws.Cells.Item(181, 8).Value = "''626402.175.524111.03212WA.2378EBA.A000000001.00000.1.0151.2.000000.000000.957.101.0B.000322"
$ws.Cells.Item(181, 9).Value = "IDR"
$ws.Cells.Item(181, 10).Value = 1
$ws.Cells.Item(181, 11).Value = 1205000
$ws.Cells.Item(181, 12).Value = 1205000

$ws.Cells.Item(182, 1).Value = 179
$ws.Cells.Item(182, 2).Value = 626402
$ws.Cells.Item(182, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(182, 4).Value = "2024-02-06"
$ws.Cells.Item(182, 5).Value = "''00063T"
$ws.Cells.Item(182, 6).Value = "''241751303001087"
$ws.Cells.Item(182, 7).Value = "Pembayaran belanja barangberupa biaya perjadin sesuai STNo.B.309,B.356,B.729/BPPSDM.1/KP.440/I/2024 Tanggal 15,16 dan 19 Januari 2024 a.n Achmad Irfansyah, dkk"
$ws.Cells.Item(182, 8).Value = "''626402.175.524113.03212WA.2378EBA.A000000001.00000.1.0151.2.000000.000000.957.101.0B.000323"
$ws.Cells.Item(182, 9).Value = "IDR"
$ws.Cells.Item(182, 10).Value = 1
$ws.Cells.Item(182, 11).Value = 630000
$ws.Cells.Item(182, 12).Value = 630000

$ws.Cells.Item(183, 1).Value = 180
$ws.Cells.Item(183, 2).Value = 626402
$ws.Cells.Item(183, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(183, 4).Value = "2024-02-06"
$ws.Cells.Item(183, 5).Value = "''00063T"
$ws.Cells.Item(183, 6).Value = "''241751303001087"
$ws.Cells.Item(183, 7).Value = "Pembayaran belanja barangberupa biaya perjadin sesuai STNo.B.309,B.356,B.729/BPPSDM.1/KP.440/I/2024 Tanggal 15,16 dan 19 Januari 2024 a.n Achmad Irfansyah, dkk"
$ws.Cells.Item(183, 8).Value = "''626402.175.524113.03212WA.2378EBA.A000000001.00000.1.0151.2.000000.000000.957.101.0B.000324"
$ws.Cells.Item(183, 9).Value = "IDR"
$ws.Cells.Item(183, 10).Value = 1
$ws.Cells.Item(183, 11).Value = 510000
$ws.Cells.Item(183, 12).Value = 510000

$ws.Cells.Item(184, 1).Value = 181
$ws.Cells.Item(184, 2).Value = 626402
$ws.Cells.Item(184, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(184, 4).Value = "2024-02-06"
$ws.Cells.Item(184, 5).Value = "''00064T"
$ws.Cells.Item(184, 6).Value = "''241751303001088"
$ws.Cells.Item(184, 7).Value = "Pembayaran Belanja Barang Sesuai Surat Tugas Nomor:B.206,B.261,B.278,B.704/BRSDM.1/KP.440/I/2023 Tgl.9,11,12,17 Januari 2024 Perjadin an.Andriawan Doni Purnama,DKK."
$ws.Cells.Item(184, 8).Value = "''626402.175.524114.03212WA.2378EBA.A000000001.00000.1.0151.2.000000.000000.962.101.0A.000602"
$ws.Cells.Item(184, 9).Value = "IDR"
$ws.Cells.Item(184, 10).Value = 1
$ws.Cells.Item(184, 11).Value = 1700000
$ws.Cells.Item(184, 12).Value = 1700000

$ws.Cells.Item(185, 1).Value = 182
$ws.Cells.Item(185, 2).Value = 626402
$ws.Cells.Item(185, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(185, 4).Value = "2024-02-06"
$ws.Cells.Item(185, 5).Value = "''00065T"
$ws.Cells.Item(185, 6).Value = "''241751303001091"
$ws.Cells.Item(185, 7).Value = "Pembayaran Belanja Barang Sesuai Surat Tugas Nomor:B.49,B.258/BRSDM.1/KP.440/I/2024 Tgl.9,11 Januari 2024 Perjadin an.Rudi Alek Wahyudin, DKK. "
$ws.Cells.Item(185, 8).Value = "''626402.175.524111.03212WA.2378EBD.A000000001.00000.1.0151.2.000000.000000.974.101.0A.001033"
$ws.Cells.Item(185, 9).Value = "IDR"
$ws.Cells.Item(185, 10).Value = 1
$ws.Cells.Item(185, 11).Value = 5907000
$ws.Cells.Item(185, 12).Value = 5907000

$ws.Cells.Item(186, 1).Value = 183
$ws.Cells.Item(186, 2).Value = 626402
$ws.Cells.Item(186, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(186, 4).Value = "2024-02-06"
$ws.Cells.Item(186, 5).Value = "''00066T"
$ws.Cells.Item(186, 6).Value = "''241751303001092"
$ws.Cells.Item(186, 7).Value = "Pembayaran belanja barang berupa biaya perjadin sesuai ST No. B.235, B.698, B.728/BPPSDM.1/KP.440/I/2024 Tanggal 10, 17 dan 19 Januari 2024 a.n Niken Financia"
$ws.Cells.Item(186, 8).Value = "''626402.175.524111.03212WA.2378EBD.A000000001.00000.1.0151.2.000000.000000.952.101.0D.000941"
$ws.Cells.Item(186, 9).Value = "IDR"
$ws.Cells.Item(186, 10).Value = 1
$ws.Cells.Item(186, 11).Value = 3897000
$ws.Cells.Item(186, 12).Value = 3897000

$ws.Cells.Item(187, 1).Value = 184
$ws.Cells.Item(187, 2).Value = 626402
$ws.Cells.Item(187, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(187, 4).Value = "2024-02-06"
$ws.Cells.Item(187, 5).Value = "''00067T"
$ws.Cells.Item(187, 6).Value = "''241751303001093"
$ws.Cells.Item(187, 7).Value = "Pembayaran Belanja Barang Sesuai Surat Tugas Nomor:B.342/BRSDM.1/KP.440/I/2024 Tgl.16 Januari 2024 Perjadin an.Wawan Nurliansyah,DKK."
$ws.Cells.Item(187, 8).Value = "''626402.175.524111.03212WA.2378EBA.A000000001.00000.1.0151.2.000000.000000.956.101.AA.000256"
$ws.Cells.Item(187, 9).Value = "IDR"
$ws.Cells.Item(187, 10).Value = 1
$ws.Cells.Item(187, 11).Value = 5410000
$ws.Cells.Item(187, 12).Value = 5410000

$ws.Cells.Item(188, 1).Value = 185
$ws.Cells.Item(188, 2).Value = 626402
$ws.Cells.Item(188, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(188, 4).Value = "2024-02-06"
$ws.Cells.Item(188, 5).Value = "''00068T"
$ws.Cells.Item(188, 6).Value = "''241751303001094"
$ws.Cells.Item(188, 7).Value = "Pembayaran belanja barang berupa biaya perjadin sesuai ST No. B.158, B.297 /BPPSDM.1/KP.440/I/2024 Tanggal 8 dan 15 Januari 2024 a.n Nurkholis Abellian, dkk"
$ws.Cells.Item(188, 8).Value = "''626402.175.524113.03212WA.2378EBD.A000000001.00000.1.0151.2.000000.000000.953.101.0B.001515"
$ws.Cells.Item(188, 9).Value = "IDR"
$ws.Cells.Item(188, 10).Value = 1
$ws.Cells.Item(188, 11).Value = 520000
$ws.Cells.Item(188, 12).Value = 520000

$ws.Cells.Item(189, 1).Value = 186
$ws.Cells.Item(189, 2).Value = 626402
$ws.Cells.Item(189, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(189, 4).Value = "2024-02-06"
$ws.Cells.Item(189, 5).Value = "''00068T"
$ws.Cells.Item(189, 6).Value = "''241751303001094"
$ws.Cells.Item(189, 7).Value = "Pembayaran belanja barang berupa biaya perjadin sesuai ST No. B.158, B.297 /BPPSDM.1/KP.440/I/2024 Tanggal 8 dan 15 Januari 2024 a.n Nurkholis Abellian, dkk"
$ws.Cells.Item(189, 8).Value = "''626402.175.524113.03212WA.2378EBD.A000000001.00000.1.0151.2.000000.000000.953.101.0B.001516"
$ws.Cells.Item(189, 9).Value = "IDR"
$ws.Cells.Item(189, 10).Value = 1
$ws.Cells.Item(189, 11).Value = 850000
$ws.Cells.Item(189, 12).Value = 850000

$ws.Cells.Item(190, 1).Value = 187
$ws.Cells.Item(190, 2).Value = 626402
$ws.Cells.Item(190, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(190, 4).Value = "2024-02-06"
$ws.Cells.Item(190, 5).Value = "''00069T"
$ws.Cells.Item(190, 6).Value = "''241751303001095"
$ws.Cells.Item(190, 7).Value = "Pembayaran belanja barang berupa biaya perjadin sesuai ST No. B.263, B.728/BPPSDM.1/KP.440/I/2024 Tanggal 10 dan 19 Januari 2024 a.n Setyawati, dkk"
$ws.Cells.Item(190, 8).Value = "''626402.175.524111.03212WA.2378EBD.A000000001.00000.1.0151.2.000000.000000.952.101.0D.000940"
$ws.Cells.Item(190, 9).Value = "IDR"
$ws.Cells.Item(190, 10).Value = 1
$ws.Cells.Item(190, 11).Value = 29064400
$ws.Cells.Item(190, 12).Value = 29064400

$ws.Cells.Item(191, 1).Value = 188
$ws.Cells.Item(191, 2).Value = 626402
$ws.Cells.Item(191, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(191, 4).Value = "2024-02-06"
$ws.Cells.Item(191, 5).Value = "''00070T"
$ws.Cells.Item(191, 6).Value = "''241751303001086"
$ws.Cells.Item(191, 7).Value = "Pembayaran belanja barang berupa biaya perjadin sesuai ST No. B.311/BPPSDM.1/KP.440/I/2024 Tanggal 15 Januari 2024 a.n Rudi Alek Wahyudin, dkk"
$ws.Cells.Item(191, 8).Value = "''626402.175.524114.03212WA.2378EBD.A000000001.00000.1.0151.2.000000.000000.952.101.0D.000943"
$ws.Cells.Item(191, 9).Value = "IDR"
$ws.Cells.Item(191, 10).Value = 1
$ws.Cells.Item(191, 11).Value = 40120000
$ws.Cells.Item(191, 12).Value = 40120000

$ws.Cells.Item(192, 1).Value = 189
$ws.Cells.Item(192, 2).Value = 626402
$ws.Cells.Item(192, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(192, 4).Value = "2024-02-06"
$ws.Cells.Item(192, 5).Value = "''00070T"
$ws.Cells.Item(192, 6).Value = "''241751303001086"
$ws.Cells.Item(192, 7).Value = "Pembayaran belanja barang berupa biaya perjadin sesuai ST No. B.311/BPPSDM.1/KP.440/I/2024 Tanggal 15 Januari 2024 a.n Rudi Alek Wahyudin, dkk"
$ws.Cells.Item(192, 8).Value = "''626402.175.524114.03212WA.2378EBD.A000000001.00000.1.0151.2.000000.000000.952.101.0D.000944"
$ws.Cells.Item(192, 9).Value = "IDR"
$ws.Cells.Item(192, 10).Value = 1
$ws.Cells.Item(192, 11).Value = 17650000
$ws.Cells.Item(192, 12).Value = 17650000

$ws.Cells.Item(193, 1).Value = 190
$ws.Cells.Item(193, 2).Value = 626402
$ws.Cells.Item(193, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(193, 4).Value = "2024-02-06"
$ws.Cells.Item(193, 5).Value = "''00071T"
$ws.Cells.Item(193, 6).Value = "''241751301001456"
$ws.Cells.Item(193, 7).Value = "Pembayaran Belanja Barang Sesuai Kuitansi Nomor:042/NNI/I/24 Tgl.29 Januari 2024 "
$ws.Cells.Item(193, 8).Value = "''626402.175.521111.03212WA.2378EBA.A000000001.00000.1.0151.2.000000.000000.994.002.AA.000425"
$ws.Cells.Item(193, 9).Value = "IDR"
$ws.Cells.Item(193, 10).Value = 1
$ws.Cells.Item(193, 11).Value = 20000000
$ws.Cells.Item(193, 12).Value = 20000000

$ws.Cells.Item(194, 1).Value = 191
$ws.Cells.Item(194, 2).Value = 626402
$ws.Cells.Item(194, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(194, 4).Value = "2024-02-06"
$ws.Cells.Item(194, 5).Value = "''00072T"
$ws.Cells.Item(194, 6).Value = "''241751301001452"
$ws.Cells.Item(194, 7).Value = "Pembayaran Belanja Barang Sesuai Kuitansi Nomor:011/PO-Sales/Kompas.com/I/2024 tgl.29 Januari 2024"
$ws.Cells.Item(194, 8).Value = "''626402.175.522191.03212WA.2378EBA.A000000001.00000.1.0151.2.000000.000000.958.101.0A.000343"
$ws.Cells.Item(194, 9).Value = "IDR"
$ws.Cells.Item(194, 10).Value = 1
$ws.Cells.Item(194, 11).Value = 23931600
$ws.Cells.Item(194, 12).Value = 23931600

$ws.Cells.Item(195, 1).Value = 192
$ws.Cells.Item(195, 2).Value = 626402
$ws.Cells.Item(195, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(195, 4).Value = "2024-02-07"
$ws.Cells.Item(195, 5).Value = "''00073T"
$ws.Cells.Item(195, 6).Value = "''241751301001568"
$ws.Cells.Item(195, 7).Value = "Pembayaran Belanja Barang-Sesuai SPK No:1230/SPK/PPBJ/BPPSDM.5/XII/2023 Tgl.29-12-2023,BAST No:BAST.1231/PPBJ.PL/BPPSDM.5/I/2024.Tgl.31-1-2024,BAP No:1231/PPBJ.PL/BPPSDM.5/II/2024 Tgl.01-02-2024.Pembayaran Bulan Januari 2024  "
$ws.Cells.Item(195, 8).Value = "''626402.175.522191.03212WA.2378EBA.A000000001.00000.1.0151.2.000000.000000.994.002.DA.000539"
$ws.Cells.Item(195, 9).Value = "IDR"
$ws.Cells.Item(195, 10).Value = 1
$ws.Cells.Item(195, 11).Value = 28720670
$ws.Cells.Item(195, 12).Value = 28720670

$ws.Cells.Item(196, 1).Value = 193
$ws.Cells.Item(196, 2).Value = 626402
$ws.Cells.Item(196, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(196, 4).Value = "2024-02-07"
$ws.Cells.Item(196, 5).Value = "''00074T"
$ws.Cells.Item(196, 6).Value = "''241751302001394"
$ws.Cells.Item(196, 7).Value = "Pembayaran Belanja Barang-Sesuai SPK No:1220/SPK/PPBJ/BPPSDM.5/XII/2023 Tgl.29-12-2023,BAST No:BAST.1221/PPBJ.PL/BPPSDM.5/I/2024.Tgl.31-1-2024,BAP No:1221/PPBJ.PL/BPPSDM.5/II/2024 Tgl.01-02-2024.Pembayaran Bulan Januari 2024 "
$ws.Cells.Item(196, 8).Value = "''626402.175.522191.03212WA.2378EBA.A000000001.00000.1.0151.2.000000.000000.994.002.DA.000540"
$ws.Cells.Item(196, 9).Value = "IDR"
$ws.Cells.Item(196, 10).Value = 1
$ws.Cells.Item(196, 11).Value = 28720670
$ws.Cells.Item(196, 12).Value = 28720670

$ws.Cells.Item(197, 1).Value = 194
$ws.Cells.Item(197, 2).Value = 626402
$ws.Cells.Item(197, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(197, 4).Value = "2024-02-07"
$ws.Cells.Item(197, 5).Value = "''00075T"
$ws.Cells.Item(197, 6).Value = "''241751302001395"
$ws.Cells.Item(197, 7).Value = "Pembayaran Belanja Barang-Sesuai SPK No:1210/SPK/PPBJ/BPPSDM.5/XII/2023 Tgl.29-12-2023,BAST No:BAST.1211/PPBJ.PL/BPPSDM.5/I/2024.Tgl.31-1-2024,BAP No:1211/PPBJ.PL/BPPSDM.5/II/2024 Tgl.01-02-2024.Pembayaran Bulan Januari 2024 "
$ws.Cells.Item(197, 8).Value = "''626402.175.522191.03212WA.2378EBA.A000000001.00000.1.0151.2.000000.000000.994.002.DA.000540"
$ws.Cells.Item(197, 9).Value = "IDR"
$ws.Cells.Item(197, 10).Value = 1
$ws.Cells.Item(197, 11).Value = 7428768
$ws.Cells.Item(197, 12).Value = 7428768

$ws.Cells.Item(198, 1).Value = 195
$ws.Cells.Item(198, 2).Value = 626402
$ws.Cells.Item(198, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(198, 4).Value = "2024-02-12"
$ws.Cells.Item(198, 5).Value = "''00076T"
$ws.Cells.Item(198, 6).Value = "''241751302001665"
$ws.Cells.Item(198, 7).Value = "Pembayaran Belanja Barang Sesuai Kuitansi Nomor:021/KW-AH/I/2024 Tanggal 30 Januari 2024"
$ws.Cells.Item(198, 8).Value = "''626402.175.521111.03212WA.2378EBA.A000000001.00000.1.0151.2.000000.000000.994.002.AA.000428"
$ws.Cells.Item(198, 9).Value = "IDR"
$ws.Cells.Item(198, 10).Value = 1
$ws.Cells.Item(198, 11).Value = 14600000
$ws.Cells.Item(198, 12).Value = 14600000

$ws.Cells.Item(199, 1).Value = 196
$ws.Cells.Item(199, 2).Value = 626402
$ws.Cells.Item(199, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(199, 4).Value = "2024-02-07"
$ws.Cells.Item(199, 5).Value = "''00077T"
$ws.Cells.Item(199, 6).Value = "''241751303001168"
$ws.Cells.Item(199, 7).Value = "Pembayaran Belanja Barang Sesuai Surat Tugas Nomor:B.271/BPPSDM.1/KP.440/I/2024 Tgl. 12 Januari 2024 Perjadin an.Ollyvia Maria Christy, DKK."
$ws.Cells.Item(199, 8).Value = "''626402.175.524111.03212WA.2378EBD.A000000001.00000.1.0151.2.000000.000000.955.101.AA.001066"
$ws.Cells.Item(199, 9).Value = "IDR"
$ws.Cells.Item(199, 10).Value = 1
$ws.Cells.Item(199, 11).Value = 1356500
$ws.Cells.Item(199, 12).Value = 1356500

$ws.Cells.Item(200, 1).Value = 197
$ws.Cells.Item(200, 2).Value = 626402
$ws.Cells.Item(200, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(200, 4).Value = "2024-02-07"
$ws.Cells.Item(200, 5).Value = "''00078T"
$ws.Cells.Item(200, 6).Value = "''241751303001420"
$ws.Cells.Item(200, 7).Value = "Penggantian uang persediaan KKP untuk keperluan belanja barang (BPP PuslatluhKP)"
$ws.Cells.Item(200, 8).Value = "''626402.175.524111.03212DL.2375AFA.A000000001.00000.1.0151.2.000000.000000.001.051.0F.001423"
$ws.Cells.Item(200, 9).Value = "IDR"
$ws.Cells.Item(200, 10).Value = 1
$ws.Cells.Item(200, 11).Value = 8099000
$ws.Cells.Item(200, 12).Value = 8099000

$ws.Cells.Item(201, 1).Value = 198
$ws.Cells.Item(201, 2).Value = 626402
$ws.Cells.Item(201, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(201, 4).Value = "2024-02-07"
$ws.Cells.Item(201, 5).Value = "''00078T"
$ws.Cells.Item(201, 6).Value = "''241751303001420"
$ws.Cells.Item(201, 7).Value = "Penggantian uang persediaan KKP untuk keperluan belanja barang (BPP PuslatluhKP)"
$ws.Cells.Item(201, 8).Value = "''626402.175.524111.03212WA.2378EBD.A000000001.00000.1.0151.2.000000.000000.952.201.0B.001473"
$ws.Cells.Item(201, 9).Value = "IDR"
$ws.Cells.Item(201, 10).Value = 1
$ws.Cells.Item(201, 11).Value = 2054100
$ws.Cells.Item(201, 12).Value = 2054100

$ws.Cells.Item(202, 1).Value = 199
$ws.Cells.Item(202, 2).Value = 626402
$ws.Cells.Item(202, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(202, 4).Value = "2024-02-07"
$ws.Cells.Item(202, 5).Value = "''00078T"
$ws.Cells.Item(202, 6).Value = "''241751303001420"
$ws.Cells.Item(202, 7).Value = "Penggantian uang persediaan KKP untuk keperluan belanja barang (BPP PuslatluhKP)"
$ws.Cells.Item(202, 8).Value = "''626402.175.524111.03212WA.2378EBD.A000000001.00000.1.0151.2.000000.000000.952.201.0D.001478"
$ws.Cells.Item(202, 9).Value = "IDR"
$ws.Cells.Item(202, 10).Value = 1
$ws.Cells.Item(202, 11).Value = 3485000
$ws.Cells.Item(202, 12).Value = 3485000

$ws.Cells.Item(203, 1).Value = 200
$ws.Cells.Item(203, 2).Value = 626402
$ws.Cells.Item(203, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(203, 4).Value = "2024-02-07"
$ws.Cells.Item(203, 5).Value = "''00078T"
$ws.Cells.Item(203, 6).Value = "''241751303001420"
$ws.Cells.Item(203, 7).Value = "Penggantian uang persediaan KKP untuk keperluan belanja barang (BPP PuslatluhKP)"
$ws.Cells.Item(203, 8).Value = "''626402.175.524111.03212WA.2378EBD.A000000001.00000.1.0151.2.000000.000000.952.201.0D.001479"
$ws.Cells.Item(203, 9).Value = "IDR"
$ws.Cells.Item(203, 10).Value = 1
$ws.Cells.Item(203, 11).Value = 1560000
$ws.Cells.Item(203, 12).Value = 1560000

$ws.Cells.Item(204, 1).Value = 201
$ws.Cells.Item(204, 2).Value = 626402
$ws.Cells.Item(204, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(204, 4).Value = "2024-02-12"
$ws.Cells.Item(204, 5).Value = "''00083T"
$ws.Cells.Item(204, 6).Value = "''241751303001234"
$ws.Cells.Item(204, 7).Value = "Pembayaran Belanja Barang berupa biaya perjadin sesuai ST No. B.135,B.277,B.332,B.688,B.712/BPPSDM.1/KP.440/I/2024 Tanggal 4,12,16,17, dan 18 Januari 2024 a.n Hari Purwanto, dkk"
$ws.Cells.Item(204, 8).Value = "''626402.175.524111.03212WA.2378EBA.A000000001.00000.1.0151.2.000000.000000.960.101.0A.000788"
$ws.Cells.Item(204, 9).Value = "IDR"
$ws.Cells.Item(204, 10).Value = 1
$ws.Cells.Item(204, 11).Value = 2935600
$ws.Cells.Item(204, 12).Value = 2935600

$ws.Cells.Item(205, 1).Value = 202
$ws.Cells.Item(205, 2).Value = 626402
$ws.Cells.Item(205, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(205, 4).Value = "2024-02-12"
$ws.Cells.Item(205, 5).Value = "''00083T"
$ws.Cells.Item(205, 6).Value = "''241751303001234"
$ws.Cells.Item(205, 7).Value = "Pembayaran Belanja Barang berupa biaya perjadin sesuai ST No. B.135,B.277,B.332,B.688,B.712/BPPSDM.1/KP.440/I/2024 Tanggal 4,12,16,17, dan 18 Januari 2024 a.n Hari Purwanto, dkk"
$ws.Cells.Item(205, 8).Value = "''626402.175.524111.03212WA.2378EBA.A000000001.00000.1.0151.2.000000.000000.960.102.0A.000805"
$ws.Cells.Item(205, 9).Value = "IDR"
$ws.Cells.Item(205, 10).Value = 1
$ws.Cells.Item(205, 11).Value = 108631940
$ws.Cells.Item(205, 12).Value = 108631940

$ws.Cells.Item(206, 1).Value = 203
$ws.Cells.Item(206, 2).Value = 626402
$ws.Cells.Item(206, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(206, 4).Value = "2024-02-12"
$ws.Cells.Item(206, 5).Value = "''00083T"
$ws.Cells.Item(206, 6).Value = "''241751303001234"
$ws.Cells.Item(206, 7).Value = "Pembayaran Belanja Barang berupa biaya perjadin sesuai ST No. B.135,B.277,B.332,B.688,B.712/BPPSDM.1/KP.440/I/2024 Tanggal 4,12,16,17, dan 18 Januari 2024 a.n Hari Purwanto, dkk"
$ws.Cells.Item(206, 8).Value = "''626402.175.524114.03212WA.2378EBC.A000000001.00000.1.0151.2.000000.000000.954.103.0A.000860"
$ws.Cells.Item(206, 9).Value = "IDR"
$ws.Cells.Item(206, 10).Value = 1
$ws.Cells.Item(206, 11).Value = 150000
$ws.Cells.Item(206, 12).Value = 150000

$ws.Cells.Item(207, 1).Value = 204
$ws.Cells.Item(207, 2).Value = 626402
$ws.Cells.Item(207, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(207, 4).Value = "2024-02-12"
$ws.Cells.Item(207, 5).Value = "''00085T"
$ws.Cells.Item(207, 6).Value = "''241751303001528"
$ws.Cells.Item(207, 7).Value = "Pembayaran tunjangan kinerja susulan Pusluh bulan Januari tahun 2024 untuk 33 Pegawai/Anggota Polri/Prajurit TNI."
$ws.Cells.Item(207, 8).Value = "''626402.175.512411.03212WA.2378EBA.A000000001.00000.1.0151.2.000000.000000.994.001.DA.000410"
$ws.Cells.Item(207, 9).Value = "IDR"
$ws.Cells.Item(207, 10).Value = 1
$ws.Cells.Item(207, 11).Value = 198871146
$ws.Cells.Item(207, 12).Value = 198871146

$ws.Cells.Item(208, 1).Value = 205
$ws.Cells.Item(208, 2).Value = 626402
$ws.Cells.Item(208, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(208, 4).Value = "2024-02-12"
$ws.Cells.Item(208, 5).Value = "''00086T"
$ws.Cells.Item(208, 6).Value = "''241751303001563"
$ws.Cells.Item(208, 7).Value = "Pembayaran tunjangan kinerja susulan Pegawai Sekretariat BRSDM bulan Januari tahun 2024 untuk 74 Pegawai/Anggota Polri/Prajurit TNI."
$ws.Cells.Item(208, 8).Value = "''626402.175.512411.03212WA.2378EBA.A000000001.00000.1.0151.2.000000.000000.994.001.AA.000399"
$ws.Cells.Item(208, 9).Value = "IDR"
$ws.Cells.Item(208, 10).Value = 1
$ws.Cells.Item(208, 11).Value = 412680289
$ws.Cells.Item(208, 12).Value = 412680289

$ws.Cells.Item(209, 1).Value = 206
$ws.Cells.Item(209, 2).Value = 626402
$ws.Cells.Item(209, 3).Value = "SEKRETARIAT BADAN RISET DAN SUMBERDAYA MANUSIA KELAUTAN DAN PERIKANAN"
$ws.Cells.Item(209, 4).Value = "2024-02-12"
$ws.Cells.Item(209, 5).Value = "''00087T"
$ws.Cells.Item(209, 6).Value = "''241751303001562"
$ws.Cells.Item(209, 7).Value = "Pembayaran tunjangan kinerja susulan bulan Januari tahun 2024 untuk 40 Pegawai/Anggota Polri/Prajurit TNI (Pusdik KP)"
$ws.Cells.Item(209, 8).Value = "''626402.175.512411.03212WA.2378EBA.A000000001.00000.1.0151.2.000000.000000.994.001.GA.000422"
$ws.Cells.Item(209, 9).Value = "IDR"
$ws.Cells.Item(209, 10).Value = 1
$ws.Cells.Item(209, 11).Value = 203211317
$ws.Cells.Item(209, 12).Value = 203211317
